$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$cellRef,
        [string]$value
    )
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "315.58"
Set-TextValue "E2" "3.32%"
Set-TextValue "D3" "35.53"
Set-TextValue "E3" "-1.42%"
Set-TextValue "D4" "5.126"
Set-TextValue "E4" "0.50%"
Set-TextValue "D5" "0.08114"
Set-TextValue "E5" "2.86%"
Set-TextValue "E6" "-0.51%"
Set-TextValue "D7" "8.003"
Set-TextValue "E7" "0.90%"
Set-TextValue "D8" "4.150"
Set-TextValue "E8" "1.08%"
Set-TextValue "D9" "0.9266"
Set-TextValue "E9" "0.46%"
Set-TextValue "E10" "4.15%"
Set-TextValue "D11" "0.1878"
Set-TextValue "E11" "0.94%"
Set-TextValue "D12" "0.09233"
Set-TextValue "E12" "6.47%"
Set-TextValue "D13" "0.03606"
Set-TextValue "E13" "1.38%"
Set-TextValue "D14" "0.09903"
Set-TextValue "E14" "-0.39%"
Set-TextValue "D15" "0.001441"
Set-TextValue "E15" "0.86%"
Set-TextValue "D16" "0.005712"
Set-TextValue "E16" "1.42%"
Set-TextValue "D17" "3.473"
Set-TextValue "E17" "0.44%"
Set-TextValue "E18" "2.06%"
Set-TextValue "D19" "0.3421"
Set-TextValue "E19" "0.76%"
Set-TextValue "D20" "0.1331"
Set-TextValue "E20" "1.03%"
Set-TextValue "D21" "5.180"
Set-TextValue "E21" "0.47%"
Set-TextValue "E22" "11.24%"
Set-TextValue "D23" "0.04565"
Set-TextValue "E23" "0.09%"
Set-TextValue "E24" "1.07%"
Set-TextValue "E25" "-6.85%"
Set-TextValue "D27" "0.0004506"
Set-TextValue "E27" "-4.97%"
Set-TextValue "D39" "0.01959"
Set-TextValue "E39" "5.44%"
Set-TextValue "D40" "0.04872"
Set-TextValue "E40" "2.00%"
Set-TextValue "D41" "0.007721"
Set-TextValue "E41" "2.43%"
Set-TextValue "D42" "0.1390"
Set-TextValue "E42" "-0.68%"
Set-TextValue "D43" "0.007838"
Set-TextValue "E43" "1.31%"
Set-TextValue "D44" "0.002105"
Set-TextValue "E44" "-5.45%"
Set-TextValue "D45" "0.01164"
Set-TextValue "E45" "5.50%"
Set-TextValue "D46" "0.00006537"
Set-TextValue "E46" "3.46%"
Set-TextValue "E47" "0.28%"
Set-TextValue "D48" "39.22"
Set-TextValue "E48" "-17.36%"
Set-TextValue "E49" "-14.74%"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.28%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.28%"

Write-Host "Updated 68 cells across 29 rows (price/volume refresh)"
